# Applies the "2 Jan 2022" commit:
#  - LoginTest sheet: update the username/password test data
#  - SearchTest sheet: replace the AccountName/RunMode columns with a new
#    SearchKey/Location/RunMode layout and data
#  - Selections updated on LoginTest and SearchTest tabs to match the
#    author's final cursor position

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Test_Suite")
$ws2 = $wb.Worksheets.Item("LoginTest")
$ws3 = $wb.Worksheets.Item("SearchTest")

# --- SearchTest: new row-2 values first (keeps shared-string insertion
#     order aligned with the original authoring order: SearchTest data,
#     then LoginTest credentials, then the SearchTest header row). ---
$ws3.Range("A2").Value = "Automation Testing"
$ws3.Range("B2").Value = "Pune"

# --- LoginTest: replace the old credentials with the new ones ---
$ws2.Range("A2").Value = "mayuresh.ahirrao@gmail.com"
$ws2.Range("B2").Value = "Mayur31885"

# --- SearchTest: new header row + carry RunMode to the new 3rd column ---
$ws3.Range("A1").Value = "SearchKey"
$ws3.Range("B1").Value = "Location"
$ws3.Range("C1").Value = "RunMode"
$ws3.Range("C2").Value = "Y"

# --- Column widths on SearchTest: widen column A for the longer header/
#     data text and size the new column B ---
$ws3.Columns.Item(1).ColumnWidth = 16.5
$ws3.Columns.Item(2).ColumnWidth = 11.666666666666666

# --- Selections: move the cursor on LoginTest to B1 and on SearchTest to
#     G14, matching the saved workbook state ---
$ws2.Activate() | Out-Null
$ws2.Range("B1").Select() | Out-Null

$ws3.Activate() | Out-Null
$ws3.Range("G14").Select() | Out-Null
